# Make headings in Manifest first-class headings.
#
# The front-matter "Manifest" section uses manually-formatted paragraphs
# (direct spacing/outlineLvl/font formatting) instead of real paragraph
# styles. Convert those heading paragraphs to use the built-in Heading 2 /
# Heading 3 styles (which already exist in this document's styles part and
# are used later for "Introduction", etc.), and drop the ad-hoc
# "Table of Contents" heading paragraph that duplicated the real ToC field.

$d = $word.ActiveDocument

# Map of heading text -> style to apply. Each of these strings is unique
# within the document, so Find can locate the right paragraph reliably.
$headings = @(
    @{ Text = "Manifest"; Style = "Heading 2" },
    @{ Text = "Status"; Style = "Heading 3" },
    @{ Text = "Publisher and License"; Style = "Heading 3" },
    @{ Text = "Cover Art"; Style = "Heading 3" },
    @{ Text = "Statement of Nonaffiliation"; Style = "Heading 3" },
    @{ Text = "Attribution"; Style = "Heading 3" },
    @{ Text = "Formats"; Style = "Heading 3" }
)

foreach ($h in $headings) {
    $found = $false
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $text = $p.Range.Text.TrimEnd("`r", "`n")
        if ($text -eq $h.Text) {
            $p.Style = $h.Style
            $found = $true
            break
        }
    }
}

# Remove the standalone "Table of Contents" heading paragraph entirely
# (it sat just before the "[ToC]" field placeholder and is no longer
# wanted as a separate heading).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r", "`n")
    if ($text -eq "Table of Contents") {
        $p.Range.Delete()
        break
    }
}
